$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row above row 30 (item "#24" VASTAFLAM ...), which pushes the
# existing rows 30-33 down to 31-34. The engine's native Insert() leaves the
# newly created row with freshly-synthesised styles (and no row height / no
# merged cells), so we patch those up explicitly afterwards by copying the
# formatting back from row 31 (which now holds what used to be row 30's
# content+style, an identical style pattern to what the brand-new row needs).
# ---------------------------------------------------------------------------
$ws.Rows("30:30").Insert()

# --- fill in the new row 30 values -----------------------------------------
# Column A / B use the plain "General" numeric style - a real number.
$ws.Range("A30").Value = 24

# Columns that are already formatted as Text (numFmtId 49) keep their text
# verbatim when assigned a string.
$ws.Range("C30").Value = "VASTAFLAM 50MG 20 SUGAR COATED TAB."
$ws.Range("H30").Value = "1:1"
$ws.Range("N30").Value = "36.00"

# Columns L and P carry a *numeric* display format (165 / 2 respectively) in
# this sheet, so a plain string assignment would silently be reinterpreted
# as a number (and trailing zeros would be lost). Forcing text with a
# leading apostrophe keeps them as genuine text cells, matching the source
# file (every data cell in this table - including these - is stored as a
# shared string).
$ws.Range("L30").Value = "'1"
$ws.Range("P30").Value = "'18.0000"
$ws.Range("Q30").Value = "0:1"

# --- restore formatting (style ids, merges, row height) on the new row -----
$ws.Range("A31:Q31").Copy()
$ws.Range("A30:Q30").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A30:B30").Merge()
$ws.Range("C30:G30").Merge()
$ws.Range("H30:K30").Merge()
$ws.Range("L30:M30").Merge()
$ws.Range("N30:O30").Merge()

$ws.Rows("30:30").RowHeight = 24.75
$ws.Rows("31:31").RowHeight = 25.5
$ws.Rows("32:32").RowHeight = 25.5
$ws.Rows("33:33").RowHeight = 24.75
$ws.Rows("34:34").RowHeight = 16.5

# --- renumber the "#" column for the rows that shifted down -----------------
$ws.Range("A31").Value = 25
$ws.Range("A32").Value = 26

# --- update the grand total (sum grew by the new line's sale price) --------
$ws.Range("P33").Value = 2192.95

# --- refresh the footer timestamp text --------------------------------------
$ws.Range("A34").Value = "Sunday, 28 September, 2025 1:28 PM"
